$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 68-71: new error codes 10066-10069 (Code column only)
for ($i = 0; $i -lt 4; $i++) {
    $row = 68 + $i
    $code = 10066 + $i
    $ws.Cells.Item($row, 1).Value = $code
}

# Rows 72-75: new "account" messages with full Code / String / Level data
$accountMessages = @(
    "message_10070_account_record_created_successfully",
    "message_10071_account_record_updated_successfully",
    "message_10072_account_record_deleted_successfully",
    "message_10073_account_multiple_records_deleted_successfully"
)

for ($i = 0; $i -lt $accountMessages.Length; $i++) {
    $row = 72 + $i
    $code = 10070 + $i
    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).Value = $accountMessages[$i]
    $ws.Cells.Item($row, 4).Value = "Success"
}

# Rows 76-119: new error codes 10074-10117 (Code column only)
for ($i = 0; $i -lt 44; $i++) {
    $row = 76 + $i
    $code = 10074 + $i
    $ws.Cells.Item($row, 1).Value = $code
}

# Restore view state: scroll position and active selection
$null = $ws.Range("B75").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$win.ScrollColumn = 1
